$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value2 = 0.04380560484234053
$ws.Cells.Item(2, 4).Value2 = 0.0368424112575596
$ws.Cells.Item(2, 5).Value2 = 0.08203861909023047
$ws.Cells.Item(2, 6).Value2 = 3.034143507635491
$ws.Cells.Item(2, 7).Value2 = 0.002560527247941622
$ws.Cells.Item(2, 9).Value2 = 1.846202781391767
$ws.Cells.Item(2, 10).Value2 = 0.1463079256554352
$ws.Cells.Item(2, 11).Value2 = 2.393862358430454
$ws.Cells.Item(2, 13).Value2 = 0.620270402485346
$ws.Cells.Item(2, 14).Value2 = 2.061698343930409
$ws.Cells.Item(3, 3).Value2 = 0.04367336812962108
$ws.Cells.Item(3, 4).Value2 = 0.03657573519980417
$ws.Cells.Item(3, 5).Value2 = 0.0820351312019767
$ws.Cells.Item(3, 6).Value2 = 3.020767909071111
$ws.Cells.Item(3, 7).Value2 = 0.002565767827764623
$ws.Cells.Item(3, 9).Value2 = 1.833918239254288
$ws.Cells.Item(3, 10).Value2 = 0.1466424809545543
$ws.Cells.Item(3, 11).Value2 = 2.259194455657337
$ws.Cells.Item(3, 13).Value2 = 0.5965247844002235
$ws.Cells.Item(3, 14).Value2 = 2.085840915183557
$ws.Cells.Item(4, 3).Value2 = 0.04359181736491635
$ws.Cells.Item(4, 4).Value2 = 0.03641510557461736
$ws.Cells.Item(4, 5).Value2 = 0.08206029996370923
$ws.Cells.Item(4, 6).Value2 = 3.014295659757281
$ws.Cells.Item(4, 7).Value2 = 0.0025691545258309
$ws.Cells.Item(4, 9).Value2 = 1.827504882682007
$ws.Cells.Item(4, 10).Value2 = 0.1469108138165609
$ws.Cells.Item(4, 11).Value2 = 2.177809740164605
$ws.Cells.Item(4, 13).Value2 = 0.5822841947959887
$ws.Cells.Item(4, 14).Value2 = 2.101392832557107
$ws.Cells.Item(5, 3).Value2 = 0.0435585008927184
$ws.Cells.Item(5, 4).Value2 = 0.03635043728935194
$ws.Cells.Item(5, 5).Value2 = 0.08207743417872848
$ws.Cells.Item(5, 6).Value2 = 3.012094564977659
$ws.Cells.Item(5, 7).Value2 = 0.002570577264358599
$ws.Cells.Item(5, 9).Value2 = 1.825174394589808
$ws.Cells.Item(5, 10).Value2 = 0.1470359692080478
$ws.Cells.Item(5, 11).Value2 = 2.144971016823149
$ws.Cells.Item(5, 13).Value2 = 0.5765662018672089
$ws.Cells.Item(5, 14).Value2 = 2.107913233278406
$ws.Cells.Item(6, 3).Value2 = 0.04355296382685125
$ws.Cells.Item(6, 4).Value2 = 0.03633974708390753
$ws.Cells.Item(6, 5).Value2 = 0.08208069500902049
$ws.Cells.Item(6, 6).Value2 = 3.011755396281458
$ws.Cells.Item(6, 7).Value2 = 0.002570816088118921
$ws.Cells.Item(6, 9).Value2 = 1.824804482999241
$ws.Cells.Item(6, 10).Value2 = 0.147057705540135
$ws.Cells.Item(6, 11).Value2 = 2.139537840751814
$ws.Cells.Item(6, 13).Value2 = 0.5756218726108742
$ws.Cells.Item(6, 14).Value2 = 2.109006978007613
$ws.Cells.Item(7, 3).Value2 = 0.04359136837954836
$ws.Cells.Item(7, 4).Value2 = 0.03641423022802925
$ws.Cells.Item(7, 5).Value2 = 0.08206050318001878
$ws.Cells.Item(7, 6).Value2 = 3.014264209723251
$ws.Cells.Item(7, 7).Value2 = 0.002569173540368493
$ws.Cells.Item(7, 9).Value2 = 1.827472308342692
$ws.Cells.Item(7, 10).Value2 = 0.146912437718111
$ws.Cells.Item(7, 11).Value2 = 2.177365546378212
$ws.Cells.Item(7, 13).Value2 = 0.5822067353313187
$ws.Cells.Item(7, 14).Value2 = 2.101480028902763
$ws.Cells.Item(8, 3).Value2 = 0.04376008713638058
$ws.Cells.Item(8, 4).Value2 = 0.03674982001905747
$ws.Cells.Item(8, 5).Value2 = 0.08203175419178166
$ws.Cells.Item(8, 6).Value2 = 3.029169611190341
$ws.Cells.Item(8, 7).Value2 = 0.002562299219839076
$ws.Cells.Item(8, 9).Value2 = 1.841732055853399
$ws.Cells.Item(8, 10).Value2 = 0.1464102143433372
$ws.Cells.Item(8, 11).Value2 = 2.347157890788822
$ws.Cells.Item(8, 13).Value2 = 0.6120123806218416
$ws.Cells.Item(8, 14).Value2 = 2.06987139052676
$ws.Cells.Item(9, 3).Value2 = 0.04408787541780868
$ws.Cells.Item(9, 4).Value2 = 0.03743232008788056
$ws.Cells.Item(9, 5).Value2 = 0.0821916953390982
$ws.Cells.Item(9, 6).Value2 = 3.072270686115402
$ws.Cells.Item(9, 7).Value2 = 0.002550152680334637
$ws.Cells.Item(9, 9).Value2 = 1.878705633400145
$ws.Cells.Item(9, 10).Value2 = 0.1459252631536181
$ws.Cells.Item(9, 11).Value2 = 2.69052430987665
$ws.Cells.Item(9, 13).Value2 = 0.6731648271711634
$ws.Cells.Item(9, 14).Value2 = 2.013676417020758
$ws.Cells.Item(10, 3).Value2 = 0.04432654282222259
$ws.Cells.Item(10, 4).Value2 = 0.03794834967643723
$ws.Cells.Item(10, 5).Value2 = 0.08244068356011702
$ws.Cells.Item(10, 6).Value2 = 3.112488384879299
$ws.Cells.Item(10, 7).Value2 = 0.002542032532882879
$ws.Cells.Item(10, 9).Value2 = 1.91143699492028
$ws.Cells.Item(10, 10).Value2 = 0.1458748981876639
$ws.Cells.Item(10, 11).Value2 = 2.949271374237014
$ws.Cells.Item(10, 13).Value2 = 0.7197633882394285
$ws.Cells.Item(10, 14).Value2 = 1.975932995520012
$ws.Cells.Item(11, 3).Value2 = 0.04443459300642516
$ws.Cells.Item(11, 4).Value2 = 0.03818622444869746
$ws.Cells.Item(11, 5).Value2 = 0.08258243258678988
$ws.Cells.Item(11, 6).Value2 = 3.132662625056412
$ws.Cells.Item(11, 7).Value2 = 0.002538511040375774
$ws.Cells.Item(11, 9).Value2 = 1.927552693280006
$ws.Cells.Item(11, 10).Value2 = 0.1459186903556713
$ws.Cells.Item(11, 11).Value2 = 3.068419432689041
$ws.Cells.Item(11, 13).Value2 = 0.7413300245570298
$ws.Cells.Item(11, 14).Value2 = 1.959534873602969
$ws.Cells.Item(12, 3).Value2 = 0.04447542885504063
$ws.Cells.Item(12, 4).Value2 = 0.0382767465473961
$ws.Cells.Item(12, 5).Value2 = 0.08264019767670305
$ws.Cells.Item(12, 6).Value2 = 3.140573901002853
$ws.Cells.Item(12, 7).Value2 = 0.002537202182115352
$ws.Cells.Item(12, 9).Value2 = 1.933832851737435
$ws.Cells.Item(12, 10).Value2 = 0.1459448859064949
$ws.Cells.Item(12, 11).Value2 = 3.11374713089657
$ws.Cells.Item(12, 13).Value2 = 0.7495500818847347
$ws.Cells.Item(12, 14).Value2 = 1.953436593648568
$ws.Cells.Item(13, 3).Value2 = 0.04446663778657722
$ws.Cells.Item(13, 4).Value2 = 0.03825723133591552
$ws.Cells.Item(13, 5).Value2 = 0.08262757524960485
$ws.Cells.Item(13, 6).Value2 = 3.138857956316343
$ws.Cells.Item(13, 7).Value2 = 0.002537482973876667
$ws.Cells.Item(13, 9).Value2 = 1.932472394038086
$ws.Cells.Item(13, 10).Value2 = 0.1459388163774733
$ws.Cells.Item(13, 11).Value2 = 3.103975684915326
$ws.Cells.Item(13, 13).Value2 = 0.7477773740032774
$ws.Cells.Item(13, 14).Value2 = 1.954745009556969
$ws.Cells.Item(14, 3).Value2 = 0.04443795423896191
$ws.Cells.Item(14, 4).Value2 = 0.03819366288966464
$ws.Cells.Item(14, 5).Value2 = 0.08258710307464945
$ws.Cells.Item(14, 6).Value2 = 3.133308035783642
$ws.Cells.Item(14, 7).Value2 = 0.002538402866384045
$ws.Cells.Item(14, 9).Value2 = 1.928065801586271
$ws.Cells.Item(14, 10).Value2 = 0.1459206527123698
$ws.Cells.Item(14, 11).Value2 = 3.072144374258471
$ws.Cells.Item(14, 13).Value2 = 0.7420052253196303
$ws.Cells.Item(14, 14).Value2 = 1.959030930927945
$ws.Cells.Item(15, 3).Value2 = 0.04442037410446531
$ws.Cells.Item(15, 4).Value2 = 0.03815478304031572
$ws.Cells.Item(15, 5).Value2 = 0.08256284481911358
$ws.Cells.Item(15, 6).Value2 = 3.129943982181146
$ws.Cells.Item(15, 7).Value2 = 0.002538969533988578
$ws.Cells.Item(15, 9).Value2 = 1.925389786961418
$ws.Cells.Item(15, 10).Value2 = 0.1459107793820422
$ws.Cells.Item(15, 11).Value2 = 3.052674028098295
$ws.Cells.Item(15, 13).Value2 = 0.7384765576906602
$ws.Cells.Item(15, 14).Value2 = 1.961670691928362
$ws.Cells.Item(16, 3).Value2 = 0.04431947053389962
$ws.Cells.Item(16, 4).Value2 = 0.03793286646568816
$ws.Cells.Item(16, 5).Value2 = 0.0824319923881891
$ws.Cells.Item(16, 6).Value2 = 3.111207891494473
$ws.Cells.Item(16, 7).Value2 = 0.002542266128383389
$ws.Cells.Item(16, 9).Value2 = 1.910408568515606
$ws.Cells.Item(16, 10).Value2 = 0.1458733802510324
$ws.Cells.Item(16, 11).Value2 = 2.941513922770241
$ws.Cells.Item(16, 13).Value2 = 0.7183613958920745
$ws.Cells.Item(16, 14).Value2 = 1.977020209616693
$ws.Cells.Item(17, 3).Value2 = 0.04425743209608157
$ws.Cells.Item(17, 4).Value2 = 0.03779752517658963
$ws.Cells.Item(17, 5).Value2 = 0.08235900806318952
$ws.Cells.Item(17, 6).Value2 = 3.10019624634856
$ws.Cells.Item(17, 7).Value2 = 0.002544332540950637
$ws.Cells.Item(17, 9).Value2 = 1.901532935473242
$ws.Cells.Item(17, 10).Value2 = 0.1458675358588835
$ws.Cells.Item(17, 11).Value2 = 2.873691333100965
$ws.Cells.Item(17, 13).Value2 = 0.7061159860058126
$ws.Cells.Item(17, 14).Value2 = 1.986634525096392
$ws.Cells.Item(18, 3).Value2 = 0.04422170052750474
$ws.Cells.Item(18, 4).Value2 = 0.03771997542854777
$ws.Cells.Item(18, 5).Value2 = 0.08231971083136713
$ws.Cells.Item(18, 6).Value2 = 3.094039407415423
$ws.Cells.Item(18, 7).Value2 = 0.002545537322617922
$ws.Cells.Item(18, 9).Value2 = 1.896543234371208
$ws.Cells.Item(18, 10).Value2 = 0.1458704514982685
$ws.Cells.Item(18, 11).Value2 = 2.834817279388972
$ws.Cells.Item(18, 13).Value2 = 0.6991074418642711
$ws.Cells.Item(18, 14).Value2 = 1.992237023454383
$ws.Cells.Item(19, 3).Value2 = 0.04420959423834603
$ws.Cells.Item(19, 4).Value2 = 0.03769376926962309
$ws.Cells.Item(19, 5).Value2 = 0.08230686624911598
$ws.Cells.Item(19, 6).Value2 = 3.091985121394629
$ws.Cells.Item(19, 7).Value2 = 0.002545948033591805
$ws.Cells.Item(19, 9).Value2 = 1.894873579049488
$ws.Cells.Item(19, 10).Value2 = 0.1458725161797005
$ws.Cells.Item(19, 11).Value2 = 2.821678456877692
$ws.Cells.Item(19, 13).Value2 = 0.6967404210519561
$ws.Cells.Item(19, 14).Value2 = 1.994146390373338
$ws.Cells.Item(20, 3).Value2 = 0.0442640412740154
$ws.Cells.Item(20, 4).Value2 = 0.03781190200362872
$ws.Cells.Item(20, 5).Value2 = 0.08236649991113865
$ws.Cells.Item(20, 6).Value2 = 3.101350148945926
$ws.Cells.Item(20, 7).Value2 = 0.002544110888574592
$ws.Cells.Item(20, 9).Value2 = 1.902465817159609
$ws.Cells.Item(20, 10).Value2 = 0.1458675081936391
$ws.Cells.Item(20, 11).Value2 = 2.880897108121019
$ws.Cells.Item(20, 13).Value2 = 0.7074159387709926
$ws.Cells.Item(20, 14).Value2 = 1.985603548385921
$ws.Cells.Item(21, 3).Value2 = 0.04444638151746716
$ws.Cells.Item(21, 4).Value2 = 0.03821232247706163
$ws.Cells.Item(21, 5).Value2 = 0.08259887985421699
$ws.Cells.Item(21, 6).Value2 = 3.134930795008387
$ws.Cells.Item(21, 7).Value2 = 0.002538132003726322
$ws.Cells.Item(21, 9).Value2 = 1.929355297914981
$ws.Cells.Item(21, 10).Value2 = 0.1459257267764826
$ws.Cells.Item(21, 11).Value2 = 3.08148832011517
$ws.Cells.Item(21, 13).Value2 = 0.7436991985910311
$ws.Cells.Item(21, 14).Value2 = 1.957769027272132
$ws.Cells.Item(22, 3).Value2 = 0.0445650802447517
$ws.Cells.Item(22, 4).Value2 = 0.03847660630241734
$ws.Cells.Item(22, 5).Value2 = 0.08277457693106882
$ws.Cells.Item(22, 6).Value2 = 3.158462261169205
$ws.Cells.Item(22, 7).Value2 = 0.002534368099580357
$ws.Cells.Item(22, 9).Value2 = 1.947964259548357
$ws.Cells.Item(22, 10).Value2 = 0.1460198144308293
$ws.Cells.Item(22, 11).Value2 = 3.213804693207464
$ws.Cells.Item(22, 13).Value2 = 0.7677228392977042
$ws.Cells.Item(22, 14).Value2 = 1.940226677092163
$ws.Cells.Item(23, 3).Value2 = 0.04450177339909089
$ws.Cells.Item(23, 4).Value2 = 0.03833531835437043
$ws.Cells.Item(23, 5).Value2 = 0.08267862683662131
$ws.Cells.Item(23, 6).Value2 = 3.145757574182653
$ws.Cells.Item(23, 7).Value2 = 0.002536363867734781
$ws.Cells.Item(23, 9).Value2 = 1.937937193935269
$ws.Cells.Item(23, 10).Value2 = 0.1459644636074984
$ws.Cells.Item(23, 11).Value2 = 3.143072976923804
$ws.Cells.Item(23, 13).Value2 = 0.7548724906946518
$ws.Cells.Item(23, 14).Value2 = 1.949529835633246
$ws.Cells.Item(24, 3).Value2 = 0.04426105346720632
$ws.Cells.Item(24, 4).Value2 = 0.03780540143024425
$ws.Cells.Item(24, 5).Value2 = 0.08236310455321316
$ws.Cells.Item(24, 6).Value2 = 3.100827928029332
$ws.Cells.Item(24, 7).Value2 = 0.00254421104540592
$ws.Cells.Item(24, 9).Value2 = 1.902043709411373
$ws.Cells.Item(24, 10).Value2 = 0.1458675011555144
$ws.Cells.Item(24, 11).Value2 = 2.877639009622953
$ws.Cells.Item(24, 13).Value2 = 0.7068281320108269
$ws.Cells.Item(24, 14).Value2 = 1.986069419065132
$ws.Cells.Item(25, 3).Value2 = 0.04399955882339768
$ws.Cells.Item(25, 4).Value2 = 0.03724510309260509
$ws.Cells.Item(25, 5).Value2 = 0.08212531033685444
$ws.Cells.Item(25, 6).Value2 = 3.059115934085753
$ws.Cells.Item(25, 7).Value2 = 0.002553296794728113
$ws.Cells.Item(25, 9).Value2 = 1.867731103773224
$ws.Cells.Item(25, 10).Value2 = 0.1460028178355302
$ws.Cells.Item(25, 11).Value2 = 2.5965071629725
$ws.Cells.Item(25, 13).Value2 = 0.6563298199032204
$ws.Cells.Item(25, 14).Value2 = 2.028257468925638

Write-Output "Updated 240 cells"